$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by Excel when writing cells.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itga4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 0.4500866666666667
$ws.Cells.Item(2, 8).Value = 1.35026
$ws.Cells.Item(2, 9).Value = 0.02628438542510526
$ws.Cells.Item(2, 10).Value = 0.02628438542510525
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.4702473333333333
$ws.Cells.Item(2, 14).Value = 1.410742
$ws.Cells.Item(2, 15).Value = 0.00903492226842282
$ws.Cells.Item(2, 16).Value = 0.00903492226842282
$ws.Cells.Item(2, 17).Value = 0.2116520547688889
$ws.Cells.Item(2, 18).Value = 1.90486849292
$ws.Cells.Item(2, 19).Value = 0.0002374773791890917
$ws.Cells.Item(2, 20).Value = 0.0002374773791890916

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itga4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 0.4500866666666667
$ws.Cells.Item(3, 8).Value = 1.35026
$ws.Cells.Item(3, 9).Value = 0.02628438542510526
$ws.Cells.Item(3, 10).Value = 0.02628438542510525
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 0.3194813333333333
$ws.Cells.Item(3, 14).Value = 0.9584440000000001
$ws.Cells.Item(3, 15).Value = 0.006138235792679485
$ws.Cells.Item(3, 16).Value = 0.006138235792679485
$ws.Cells.Item(3, 17).Value = 0.1437942883822222
$ws.Cells.Item(3, 18).Value = 1.29414859544
$ws.Cells.Item(3, 19).Value = 0.0001613397554049641
$ws.Cells.Item(3, 20).Value = 0.000161339755404964

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itga4"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.4500866666666667
$ws.Cells.Item(4, 8).Value = 1.35026
$ws.Cells.Item(4, 9).Value = 0.02628438542510526
$ws.Cells.Item(4, 10).Value = 0.02628438542510525
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 1.047307
$ws.Cells.Item(4, 14).Value = 3.141921
$ws.Cells.Item(4, 15).Value = 0.02012204358311108
$ws.Cells.Item(4, 16).Value = 0.02012204358311108
$ws.Cells.Item(4, 17).Value = 0.4713789166066667
$ws.Cells.Item(4, 18).Value = 4.24241024946
$ws.Cells.Item(4, 19).Value = 0.0005288955490792577
$ws.Cells.Item(4, 20).Value = 0.0005288955490792576

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itga4"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 0.4500866666666667
$ws.Cells.Item(5, 8).Value = 1.35026
$ws.Cells.Item(5, 9).Value = 0.02628438542510526
$ws.Cells.Item(5, 10).Value = 0.02628438542510525
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 50.21070966666667
$ws.Cells.Item(5, 14).Value = 150.632129
$ws.Cells.Item(5, 15).Value = 0.9647047983557866
$ws.Cells.Item(5, 16).Value = 0.9647047983557866
$ws.Cells.Item(5, 17).Value = 22.59917094483778
$ws.Cells.Item(5, 18).Value = 203.39253850354
$ws.Cells.Item(5, 19).Value = 0.02535667274143194
$ws.Cells.Item(5, 20).Value = 0.02535667274143194

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itga4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 15.76143266666667
$ws.Cells.Item(6, 8).Value = 47.284298
$ws.Cells.Item(6, 9).Value = 0.9204439983318276
$ws.Cells.Item(6, 10).Value = 0.9204439983318274
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 0.4702473333333333
$ws.Cells.Item(6, 14).Value = 1.410742
$ws.Cells.Item(6, 15).Value = 0.00903492226842282
$ws.Cells.Item(6, 16).Value = 0.00903492226842282
$ws.Cells.Item(6, 17).Value = 7.411771681012888
$ws.Cells.Item(6, 18).Value = 66.705945129116
$ws.Cells.Item(6, 19).Value = 0.008316139977364365
$ws.Cells.Item(6, 20).Value = 0.008316139977364365

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itga4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 15.76143266666667
$ws.Cells.Item(7, 8).Value = 47.284298
$ws.Cells.Item(7, 9).Value = 0.9204439983318276
$ws.Cells.Item(7, 10).Value = 0.9204439983318274
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 0.3194813333333333
$ws.Cells.Item(7, 14).Value = 0.9584440000000001
$ws.Cells.Item(7, 15).Value = 0.006138235792679485
$ws.Cells.Item(7, 16).Value = 0.006138235792679485
$ws.Cells.Item(7, 17).Value = 5.035483523590222
$ws.Cells.Item(7, 18).Value = 45.31935171231201
$ws.Cells.Item(7, 19).Value = 0.00564990229571744
$ws.Cells.Item(7, 20).Value = 0.005649902295717439

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itga4"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 15.76143266666667
$ws.Cells.Item(8, 8).Value = 47.284298
$ws.Cells.Item(8, 9).Value = 0.9204439983318276
$ws.Cells.Item(8, 10).Value = 0.9204439983318274
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 1.047307
$ws.Cells.Item(8, 14).Value = 3.141921
$ws.Cells.Item(8, 15).Value = 0.02012204358311108
$ws.Cells.Item(8, 16).Value = 0.02012204358311108
$ws.Cells.Item(8, 17).Value = 16.50705876182867
$ws.Cells.Item(8, 18).Value = 148.563528856458
$ws.Cells.Item(8, 19).Value = 0.01852121425024606
$ws.Cells.Item(8, 20).Value = 0.01852121425024606

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itga4"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 15.76143266666667
$ws.Cells.Item(9, 8).Value = 47.284298
$ws.Cells.Item(9, 9).Value = 0.9204439983318276
$ws.Cells.Item(9, 10).Value = 0.9204439983318274
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 50.21070966666667
$ws.Cells.Item(9, 14).Value = 150.632129
$ws.Cells.Item(9, 15).Value = 0.9647047983557866
$ws.Cells.Item(9, 16).Value = 0.9647047983557866
$ws.Cells.Item(9, 17).Value = 791.3927195567159
$ws.Cells.Item(9, 18).Value = 7122.534476010443
$ws.Cells.Item(9, 19).Value = 0.8879567418084997
$ws.Cells.Item(9, 20).Value = 0.8879567418084996

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itga4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.9026056666666666
$ws.Cells.Item(10, 8).Value = 2.707817
$ws.Cells.Item(10, 9).Value = 0.05271081546417152
$ws.Cells.Item(10, 10).Value = 0.05271081546417151
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.4702473333333333
$ws.Cells.Item(10, 14).Value = 1.410742
$ws.Cells.Item(10, 15).Value = 0.00903492226842282
$ws.Cells.Item(10, 16).Value = 0.00903492226842282
$ws.Cells.Item(10, 17).Value = 0.4244479078015555
$ws.Cells.Item(10, 18).Value = 3.820031170214
$ws.Cells.Item(10, 19).Value = 0.0004762381204239692
$ws.Cells.Item(10, 20).Value = 0.0004762381204239691

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Spon2"
$ws.Cells.Item(11, 3).Value = "Itga4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 0.9026056666666666
$ws.Cells.Item(11, 8).Value = 2.707817
$ws.Cells.Item(11, 9).Value = 0.05271081546417152
$ws.Cells.Item(11, 10).Value = 0.05271081546417151
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 0.3194813333333333
$ws.Cells.Item(11, 14).Value = 0.9584440000000001
$ws.Cells.Item(11, 15).Value = 0.006138235792679485
$ws.Cells.Item(11, 16).Value = 0.006138235792679485
$ws.Cells.Item(11, 17).Value = 0.2883656618608889
$ws.Cells.Item(11, 18).Value = 2.595290956748
$ws.Cells.Item(11, 19).Value = 0.000323551414143501
$ws.Cells.Item(11, 20).Value = 0.0003235514141435008

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Spon2"
$ws.Cells.Item(12, 3).Value = "Itga4"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 0.9026056666666666
$ws.Cells.Item(12, 8).Value = 2.707817
$ws.Cells.Item(12, 9).Value = 0.05271081546417152
$ws.Cells.Item(12, 10).Value = 0.05271081546417151
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 1.047307
$ws.Cells.Item(12, 14).Value = 3.141921
$ws.Cells.Item(12, 15).Value = 0.02012204358311108
$ws.Cells.Item(12, 16).Value = 0.02012204358311108
$ws.Cells.Item(12, 17).Value = 0.9453052329396666
$ws.Cells.Item(12, 18).Value = 8.507747096456999
$ws.Cells.Item(12, 19).Value = 0.001060649326071385
$ws.Cells.Item(12, 20).Value = 0.001060649326071385

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Spon2"
$ws.Cells.Item(13, 3).Value = "Itga4"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 0.9026056666666666
$ws.Cells.Item(13, 8).Value = 2.707817
$ws.Cells.Item(13, 9).Value = 0.05271081546417152
$ws.Cells.Item(13, 10).Value = 0.05271081546417151
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 50.21070966666667
$ws.Cells.Item(13, 14).Value = 150.632129
$ws.Cells.Item(13, 15).Value = 0.9647047983557866
$ws.Cells.Item(13, 16).Value = 0.9647047983557866
$ws.Cells.Item(13, 17).Value = 45.32047107248812
$ws.Cells.Item(13, 18).Value = 407.884239652393
$ws.Cells.Item(13, 19).Value = 0.05085037660353266
$ws.Cells.Item(13, 20).Value = 0.05085037660353265

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Spon2"
$ws.Cells.Item(14, 3).Value = "Itga4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.009603
$ws.Cells.Item(14, 8).Value = 0.028809
$ws.Cells.Item(14, 9).Value = 0.0005608007788958107
$ws.Cells.Item(14, 10).Value = 0.0005608007788958106
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 12).Value = 1.0
$ws.Cells.Item(14, 13).Value = 0.4702473333333333
$ws.Cells.Item(14, 14).Value = 1.410742
$ws.Cells.Item(14, 15).Value = 0.00903492226842282
$ws.Cells.Item(14, 16).Value = 0.00903492226842282
$ws.Cells.Item(14, 17).Value = 0.004515785142
$ws.Cells.Item(14, 18).Value = 0.040642066278
$ws.Cells.Item(14, 19).Value = 0.000005066791445394621338225213
$ws.Cells.Item(14, 20).Value = 0.000005066791445394619644159318

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Spon2"
$ws.Cells.Item(15, 3).Value = "Itga4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.009603
$ws.Cells.Item(15, 8).Value = 0.028809
$ws.Cells.Item(15, 9).Value = 0.0005608007788958107
$ws.Cells.Item(15, 10).Value = 0.0005608007788958106
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 0.3194813333333333
$ws.Cells.Item(15, 14).Value = 0.9584440000000001
$ws.Cells.Item(15, 15).Value = 0.006138235792679485
$ws.Cells.Item(15, 16).Value = 0.006138235792679485
$ws.Cells.Item(15, 17).Value = 0.003067979244
$ws.Cells.Item(15, 18).Value = 0.027611813196
$ws.Cells.Item(15, 19).Value = 0.000003442327413580798924704236
$ws.Cells.Item(15, 20).Value = 0.000003442327413580798924704236

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Spon2"
$ws.Cells.Item(16, 3).Value = "Itga4"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1.0
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.009603
$ws.Cells.Item(16, 8).Value = 0.028809
$ws.Cells.Item(16, 9).Value = 0.0005608007788958107
$ws.Cells.Item(16, 10).Value = 0.0005608007788958106
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 1.047307
$ws.Cells.Item(16, 14).Value = 3.141921
$ws.Cells.Item(16, 15).Value = 0.02012204358311108
$ws.Cells.Item(16, 16).Value = 0.02012204358311108
$ws.Cells.Item(16, 17).Value = 0.010057289121
$ws.Cells.Item(16, 18).Value = 0.090515602089
$ws.Cells.Item(16, 19).Value = 0.000011284457714384139439081113
$ws.Cells.Item(16, 20).Value = 0.000011284457714384139439081113

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Spon2"
$ws.Cells.Item(17, 3).Value = "Itga4"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.009603
$ws.Cells.Item(17, 8).Value = 0.028809
$ws.Cells.Item(17, 9).Value = 0.0005608007788958107
$ws.Cells.Item(17, 10).Value = 0.0005608007788958106
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 50.21070966666667
$ws.Cells.Item(17, 14).Value = 150.632129
$ws.Cells.Item(17, 15).Value = 0.9647047983557866
$ws.Cells.Item(17, 16).Value = 0.9647047983557866
$ws.Cells.Item(17, 17).Value = 0.4821734449290001
$ws.Cells.Item(17, 18).Value = 4.339561004361001
$ws.Cells.Item(17, 19).Value = 0.0005410072023224511
$ws.Cells.Item(17, 20).Value = 0.000541007202322451
